$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.311.04'
$ws.Range("E2").Value = '  -0.90%  '
$ws.Range("D3").Value = '1.705.71'
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '223.86'
$ws.Range("E5").Value = '  -0.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5316'
$ws.Range("E6").Value = '  -1.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2659'
$ws.Range("E8").Value = '  -0.61%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06593'
$ws.Range("E9").Value = '  -0.30%  '
$ws.Range("E10").Value = '  -4.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07659'
$ws.Range("E11").Value = '  -0.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.504'
$ws.Range("E12").Value = '  -2.40%  '
$ws.Range("D13").Value = '1.709.21'
$ws.Range("E13").Value = '  -0.69%  '
$ws.Range("D14").Value = '1.940.13'
$ws.Range("E14").Value = '  -1.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5809'
$ws.Range("E15").Value = '  -1.15%  '
$ws.Range("D16").Value = '0.0₅8166'
$ws.Range("E16").Value = '  -1.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.57'
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("D18").Value = '27.292.08'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.30'
$ws.Range("E19").Value = '  -2.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.623'
$ws.Range("E21").Value = '  -2.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.40'
$ws.Range("E22").Value = '  -2.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.986'
$ws.Range("E23").Value = '  -1.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.77'
$ws.Range("E25").Value = '  -3.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.685'
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1204'
$ws.Range("E27").Value = '  -2.45%  '
$ws.Range("E28").Value = '  -2.40%  '
$ws.Range("E29").Value = '  -2.86%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05363'
$ws.Range("E30").Value = '  -3.36%  '
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.485'
$ws.Range("E32").Value = '  -1.82%  '
$ws.Range("E33").Value = '  -1.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.645'
$ws.Range("E34").Value = '  -1.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.856'
$ws.Range("E35").Value = '  +1.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9503'
$ws.Range("E36").Value = '  -1.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.395'
$ws.Range("E37").Value = '  -2.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5844'
$ws.Range("E38").Value = '  -1.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01639'
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.807'
$ws.Range("E40").Value = '  -1.87%  '
$ws.Range("D41").Value = '1.043.80'
$ws.Range("E41").Value = '  -1.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.002'
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8416'
$ws.Range("E43").Value = '  -1.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.84'
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("D45").Value = '1.848.33'
$ws.Range("E45").Value = '  -1.06%  '
$ws.Range("D46").Value = '0.0₈117'
$ws.Range("E46").Value = '  +1.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '57.81'
$ws.Range("E47").Value = '  -2.16%  '
$ws.Range("E48").Value = '  +1.84%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.006'
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.080'
$ws.Range("E50").Value = '  -1.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05233'
$ws.Range("E51").Value = '  -0.88%  '
